$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- Remark text for each updated/new row ---
$remark2 = '1.	Updated the automation test cases in the smartsheet. Will share the updated sheet after the discussed changes are implemented by Noah.
2.	Attended daily automation discussion meeting with Noah. Discussed the changes made in Automation test cases in smartsheet.
3.	Reviewed the requirements against Christmas Cards Store Banner [#8496].
4.	Verified #8398 on mobile and desktop site. Changes are appearing fine now.'
$remark3 = '1.	Updated the automation test cases in the smartsheet. Discussed some changes with Noah. Will be ready for review, once Noah implement those changes. Also, added a Backlog section in the sheet to consolidate the scope for automation on PMall Website. Please review and provide your feedback. Also, could you please prioritize the backlog items as per the required order of implementation.
2.	Attended daily automation discussion meeting with Noah. Discussed the changes made in Automation test cases in smartsheet.
3.	Verified the ticket #7734 on prod environment. Updated results in PMall admin. Please review.
4.	Updated Test results in September release smartsheet. Please find updated sheet attached.
 '
$remark4 = '1.	Performed cross browser testing on mobile site  on iOS 12 and Android 7 around the September release tickets. No issues found.
2.	Attended daily automation discussion meeting with Noah. Discussed the changes made in Automation test cases in smartsheet. Please review the backlog section in the automation scenarios smartsheet and prioritize the scope.
3.	Discussed All production issues shared on email on 10th  and 11th October. Will create the defects against the reviewed issues in PMall admin.'
$remark5 = '1.	Created 6 new defects from #8517 - #8522 in PMall admin as per observations discussed yesterday. 
2.	Reviewed automation script code for Desktop site to understand the structure.
3.	Attended daily automation discussion meeting with Noah. Discussed the changes made in Automation test cases in smartsheet. Please review the backlog section in the automation scenarios smartsheet and prioritize the scope.
4.	Attended daily status call. Discussed tickets in September release that need to be closed or moved to October release. Will update the status of tickets in PMall admin as per discussion and update status in smartsheet too, accordingly.'
$remark6 = '1.	Created 3 new defects from #8529 - #8531 in PMall admin and added them to smartsheet as well.
2.	Closed all the tickets from September release scope as discussed yesterday. Please find the updated September release sheet attached.
3.	Worked on mobile automation scripts review in smartsheet. Will consolidate scripts with desktop ones and share the updated scripts by Monday.
4.	Attended daily automation discussion meeting with Noah. Discussed the changes made in Automation test cases in smartsheet. Also, discussed the daily automation jobs on autoqa2 server. '
$remark7 = '1.	Updated all mobile site scenarios and collaborated them with Desktop ones in the required format in automation test scenarios smartsheet. Added few comments for Noah to update. Please review and suggest.
2.	Verified few PMall pages using the Varvy  SEO inspector tool as per the assigned ticket #8394. Need to discuss it before moving on what issues to be reported in what format and also about the scope of verification.
3.	Attended daily automation discussion meeting with Noah. Discussed the changes made in Automation test cases in smartsheet.
4.	Reviewed the requirements of the October release tickets.'

# --- Extend formatting down to rows 6-12 by copying existing formatted rows ---
# Rows 6 and 7 are new data rows: copy the full A:D format from row 5 (wrap-text style in D)
$ws.Range("A5:D5").Copy()
$ws.Range("A6:D7").PasteSpecial(-4122)
$ws.Range("A6:A7").Value = "BBBY-PMall"
$ws.Range("B6:B7").Value = "Test Execution/Defect Documentation/Defect Retest/Testing"

# Rows 8-12 are blank trailer rows: copy the plain (non-wrap) A:C style across A:D
$ws.Range("A2:C2").Copy()
$ws.Range("A8:C12").PasteSpecial(-4122)
$ws.Range("A2").Copy()
$ws.Range("D8:D12").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# --- Update date (C) values for rows 2-7 first ---
$ws.Range("C2").Value = "10/14/2019"
$ws.Range("C3").Value = "10/15/2019"
$ws.Range("C4").Value = "10/16/2019"
$ws.Range("C5").Value = "10/17/2019"
$ws.Range("C6").Value = "10/18/2019"
$ws.Range("C7").Value = "10/21/2019"

# --- Then update remarks (D) values for rows 7 down to 2 (matches shared-string add order) ---
$ws.Range("D7").Value = $remark7
$ws.Range("D6").Value = $remark6
$ws.Range("D5").Value = $remark5
$ws.Range("D4").Value = $remark4
$ws.Range("D3").Value = $remark3
$ws.Range("D2").Value = $remark2

# --- Row heights to match wrapped content ---
$ws.Rows.Item(2).RowHeight = 90
$ws.Rows.Item(3).RowHeight = 180
$ws.Rows.Item(4).RowHeight = 105
$ws.Rows.Item(5).RowHeight = 135
$ws.Rows.Item(6).RowHeight = 135
$ws.Rows.Item(7).RowHeight = 135

# --- Sheet dimension grows to A1:D12 automatically; refresh used range ---
$dim = $ws.UsedRange.Address
